$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("united states")
$ws.Activate()

# Update land-grant population input values (col C) for several states
$ws.Range("C51").Value = 36239000   # New York
$ws.Range("C68").Value = 3943000    # Juniper
$ws.Range("C69").Value = 5137000    # Ontonagon
$ws.Range("C71").Value = 7821000    # Kances
$ws.Range("C74").Value = 1883000    # Alleghania
$ws.Range("C76").Value = 8304000    # Cimarron
$ws.Range("C77").Value = 3101000    # East Florida
$ws.Range("C81").Value = 2621000    # Minasota

# Update L77 formula to also add the East Florida (E259) metropolitan figure
$ws.Range("L77").Formula = "=E234*(2/3)+E259"

# New blank, number-formatted cell at M88 (extends the used range to column M)
$ws.Range("M88").Value = ""
$ws.Range("M88").NumberFormat = '_-* #,##0_-;\-* #,##0_-;_-* "-"??_-;_-@_-'

# Restore the view/selection state recorded for this sheet
$ws.Range("C52").Select()
